# Changes made on day 14 Aug 2024
#
# This workbook is used as a running log of generated test credentials.
# Every automated test run writes a fresh, unique "User Name"/"Password"
# pair into row 2 (A2/B2). Each new unique string that ever gets written
# is appended to the shared string table, while row 2 itself always shows
# only the most recently generated credentials.
#
# Replay every run that happened between the previous captured state and
# 14 Aug 2024, finishing on the credentials generated on 14 Aug 2024.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$history = @(
    @("Test12082024@221.com", "Test12082024@221"),
    @("Test12082024@421.com", "Test12082024@421"),
    @("Test12082024@590.com", "Test12082024@590"),
    @("Test12082024@584.com", "Test12082024@584"),
    @("Test13082024@458.com", "Test13082024@458"),
    @("Test13082024@681.com", "Test13082024@681"),
    @("Test13082024@266.com", "Test13082024@266"),
    @("Test13082024@361.com", "Test13082024@361"),
    @("Test13082024@939.com", "Test13082024@939"),
    @("Test13082024@301.com", "Test13082024@301"),
    @("Test13082024@977.com", "Test13082024@977"),
    @("Test13082024@942.com", "Test13082024@942"),
    @("Test13082024@973.com", "Test13082024@973"),
    @("Test13082024@848.com", "Test13082024@848"),
    @("Test13082024@763.com", "Test13082024@763"),
    @("Test13082024@945.com", "Test13082024@945"),
    @("Test13082024@604.com", "Test13082024@604"),
    @("Test13082024@366.com", "Test13082024@366"),
    @("Test13082024@269.com", "Test13082024@269"),
    @("Test13082024@0.com",   "Test13082024@0"),
    @("Test14082024@724.com", "Test14082024@724")
)

foreach ($pair in $history) {
    $ws.Range("A2").Value = $pair[0]
    $ws.Range("B2").Value = $pair[1]
}
